$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-05-12 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-05-13 Saturday", 2) | Out-Null
$d.Content.Find.Execute("98-15=83", $true, $false, $false, $false, $false, $true, 1, $false, "1+36=37", 2) | Out-Null
$d.Content.Find.Execute("30+26=56", $true, $false, $false, $false, $false, $true, 1, $false, "79-0=79", 2) | Out-Null
$d.Content.Find.Execute("43-31=12", $true, $false, $false, $false, $false, $true, 1, $false, "29-21=8", 2) | Out-Null
$d.Content.Find.Execute("8+31=39", $true, $false, $false, $false, $false, $true, 1, $false, "81-48=33", 2) | Out-Null
$d.Content.Find.Execute("6+38=44", $true, $false, $false, $false, $false, $true, 1, $false, "34+26=60", 2) | Out-Null
$d.Content.Find.Execute("9+11=20", $true, $false, $false, $false, $false, $true, 1, $false, "53-29=24", 2) | Out-Null
$d.Content.Find.Execute("76-61=15", $true, $false, $false, $false, $false, $true, 1, $false, "37+61=98", 2) | Out-Null
$d.Content.Find.Execute("69-49=20", $true, $false, $false, $false, $false, $true, 1, $false, "96-54=42", 2) | Out-Null
$d.Content.Find.Execute("34+57=91", $true, $false, $false, $false, $false, $true, 1, $false, "7+46=53", 2) | Out-Null
$d.Content.Find.Execute("35+21=56", $true, $false, $false, $false, $false, $true, 1, $false, "51-26=25", 2) | Out-Null
$d.Content.Find.Execute("35-24=11", $true, $false, $false, $false, $false, $true, 1, $false, "85-57=28", 2) | Out-Null
$d.Content.Find.Execute("68-13=55", $true, $false, $false, $false, $false, $true, 1, $false, "5+35=40", 2) | Out-Null
$d.Content.Find.Execute("1+78=79", $true, $false, $false, $false, $false, $true, 1, $false, "24+29=53", 2) | Out-Null
$d.Content.Find.Execute("69-27=42", $true, $false, $false, $false, $false, $true, 1, $false, "42-29=13", 2) | Out-Null
$d.Content.Find.Execute("36-7=29", $true, $false, $false, $false, $false, $true, 1, $false, "23+40=63", 2) | Out-Null
$d.Content.Find.Execute("39+29=68", $true, $false, $false, $false, $false, $true, 1, $false, "87-36=51", 2) | Out-Null
$d.Content.Find.Execute("88-41=47", $true, $false, $false, $false, $false, $true, 1, $false, "23+22=45", 2) | Out-Null
$d.Content.Find.Execute("73+11=84", $true, $false, $false, $false, $false, $true, 1, $false, "29+41=70", 2) | Out-Null
$d.Content.Find.Execute("88-51=37", $true, $false, $false, $false, $false, $true, 1, $false, "82-64=18", 2) | Out-Null
$d.Content.Find.Execute("40-20=20", $true, $false, $false, $false, $false, $true, 1, $false, "32-5=27", 2) | Out-Null
$d.Content.Find.Execute("89-76=13", $true, $false, $false, $false, $false, $true, 1, $false, "69+28=97", 2) | Out-Null
$d.Content.Find.Execute("35+7=42", $true, $false, $false, $false, $false, $true, 1, $false, "13-3=10", 2) | Out-Null
$d.Content.Find.Execute("92-19=73", $true, $false, $false, $false, $false, $true, 1, $false, "59+39=98", 2) | Out-Null
$d.Content.Find.Execute("11+43=54", $true, $false, $false, $false, $false, $true, 1, $false, "56+31=87", 2) | Out-Null
$d.Content.Find.Execute("5+71=76", $true, $false, $false, $false, $false, $true, 1, $false, "76-33=43", 2) | Out-Null
$d.Content.Find.Execute("96-25=71", $true, $false, $false, $false, $false, $true, 1, $false, "95-69=26", 2) | Out-Null
$d.Content.Find.Execute("36+37=73", $true, $false, $false, $false, $false, $true, 1, $false, "6+66=72", 2) | Out-Null
$d.Content.Find.Execute("96-39=57", $true, $false, $false, $false, $false, $true, 1, $false, "36-8=28", 2) | Out-Null
$d.Content.Find.Execute("26-7=19", $true, $false, $false, $false, $false, $true, 1, $false, "59-18=41", 2) | Out-Null
$d.Content.Find.Execute("0+97=97", $true, $false, $false, $false, $false, $true, 1, $false, "61-44=17", 2) | Out-Null
$d.Content.Find.Execute("99-23=76", $true, $false, $false, $false, $false, $true, 1, $false, "81-23=58", 2) | Out-Null
$d.Content.Find.Execute("50+18=68", $true, $false, $false, $false, $false, $true, 1, $false, "31+4=35", 2) | Out-Null
$d.Content.Find.Execute("26+35=61", $true, $false, $false, $false, $false, $true, 1, $false, "56-16=40", 2) | Out-Null
$d.Content.Find.Execute("60+27=87", $true, $false, $false, $false, $false, $true, 1, $false, "64+3=67", 2) | Out-Null
$d.Content.Find.Execute("60-34=26", $true, $false, $false, $false, $false, $true, 1, $false, "83-21=62", 2) | Out-Null
$d.Content.Find.Execute("35-27=8", $true, $false, $false, $false, $false, $true, 1, $false, "2+22=24", 2) | Out-Null
$d.Content.Find.Execute("11+23=34", $true, $false, $false, $false, $false, $true, 1, $false, "35+22=57", 2) | Out-Null
$d.Content.Find.Execute("31-15=16", $true, $false, $false, $false, $false, $true, 1, $false, "31-20=11", 2) | Out-Null
$d.Content.Find.Execute("94-68=26", $true, $false, $false, $false, $false, $true, 1, $false, "78-76=2", 2) | Out-Null
$d.Content.Find.Execute("81-19=62", $true, $false, $false, $false, $false, $true, 1, $false, "86+10=96", 2) | Out-Null
$d.Content.Find.Execute("57-51=6", $true, $false, $false, $false, $false, $true, 1, $false, "96-91=5", 2) | Out-Null
$d.Content.Find.Execute("95-0=95", $true, $false, $false, $false, $false, $true, 1, $false, "67-61=6", 2) | Out-Null
$d.Content.Find.Execute("98+0=98", $true, $false, $false, $false, $false, $true, 1, $false, "43+22=65", 2) | Out-Null
$d.Content.Find.Execute("33+58=91", $true, $false, $false, $false, $false, $true, 1, $false, "40-23=17", 2) | Out-Null
$d.Content.Find.Execute("41-22=19", $true, $false, $false, $false, $false, $true, 1, $false, "21+45=66", 2) | Out-Null
$d.Content.Find.Execute("99-74=25", $true, $false, $false, $false, $false, $true, 1, $false, "18+71=89", 2) | Out-Null
$d.Content.Find.Execute("36+4=40", $true, $false, $false, $false, $false, $true, 1, $false, "58-0=58", 2) | Out-Null
$d.Content.Find.Execute("6+23=29", $true, $false, $false, $false, $false, $true, 1, $false, "66+16=82", 2) | Out-Null
$d.Content.Find.Execute("81-79=2", $true, $false, $false, $false, $false, $true, 1, $false, "10+3=13", 2) | Out-Null
$d.Content.Find.Execute("47-21=26", $true, $false, $false, $false, $false, $true, 1, $false, "44+28=72", 2) | Out-Null
$d.Content.Find.Execute("67-56=11", $true, $false, $false, $false, $false, $true, 1, $false, "42-35=7", 2) | Out-Null
$d.Content.Find.Execute("40+13=53", $true, $false, $false, $false, $false, $true, 1, $false, "57-27=30", 2) | Out-Null
$d.Content.Find.Execute("68-39=29", $true, $false, $false, $false, $false, $true, 1, $false, "68-0=68", 2) | Out-Null
$d.Content.Find.Execute("57-56=1", $true, $false, $false, $false, $false, $true, 1, $false, "76-25=51", 2) | Out-Null
$d.Content.Find.Execute("84-55=29", $true, $false, $false, $false, $false, $true, 1, $false, "68-66=2", 2) | Out-Null
$d.Content.Find.Execute("6+21=27", $true, $false, $false, $false, $false, $true, 1, $false, "89-20=69", 2) | Out-Null
$d.Content.Find.Execute("85-72=13", $true, $false, $false, $false, $false, $true, 1, $false, "29+6=35", 2) | Out-Null
$d.Content.Find.Execute("18-17=1", $true, $false, $false, $false, $false, $true, 1, $false, "53-49=4", 2) | Out-Null
$d.Content.Find.Execute("21+54=75", $true, $false, $false, $false, $false, $true, 1, $false, "17+76=93", 2) | Out-Null
$d.Content.Find.Execute("38+20=58", $true, $false, $false, $false, $false, $true, 1, $false, "59+1=60", 2) | Out-Null
$d.Content.Find.Execute("52-27=25", $true, $false, $false, $false, $false, $true, 1, $false, "28+30=58", 2) | Out-Null
$d.Content.Find.Execute("44-7=37", $true, $false, $false, $false, $false, $true, 1, $false, "18-2=16", 2) | Out-Null
$d.Content.Find.Execute("18+21=39", $true, $false, $false, $false, $false, $true, 1, $false, "0+3=3", 2) | Out-Null
$d.Content.Find.Execute("15+38=53", $true, $false, $false, $false, $false, $true, 1, $false, "16+17=33", 2) | Out-Null
$d.Content.Find.Execute("22-18=4", $true, $false, $false, $false, $false, $true, 1, $false, "28-12=16", 2) | Out-Null
$d.Content.Find.Execute("3+53=56", $true, $false, $false, $false, $false, $true, 1, $false, "42+40=82", 2) | Out-Null
$d.Content.Find.Execute("94-62=32", $true, $false, $false, $false, $false, $true, 1, $false, "5+14=19", 2) | Out-Null
$d.Content.Find.Execute("30+48=78", $true, $false, $false, $false, $false, $true, 1, $false, "55-45=10", 2) | Out-Null
$d.Content.Find.Execute("40-36=4", $true, $false, $false, $false, $false, $true, 1, $false, "11+40=51", 2) | Out-Null
$d.Content.Find.Execute("9+72=81", $true, $false, $false, $false, $false, $true, 1, $false, "44+48=92", 2) | Out-Null
$d.Content.Find.Execute("93-86=7", $true, $false, $false, $false, $false, $true, 1, $false, "69-29=40", 2) | Out-Null
$d.Content.Find.Execute("38-16=22", $true, $false, $false, $false, $false, $true, 1, $false, "59-47=12", 2) | Out-Null
$d.Content.Find.Execute("36-3=33", $true, $false, $false, $false, $false, $true, 1, $false, "71+7=78", 2) | Out-Null
$d.Content.Find.Execute("39-2=37", $true, $false, $false, $false, $false, $true, 1, $false, "27+41=68", 2) | Out-Null
$d.Content.Find.Execute("88-35=53", $true, $false, $false, $false, $false, $true, 1, $false, "98-76=22", 2) | Out-Null
$d.Content.Find.Execute("40-29=11", $true, $false, $false, $false, $false, $true, 1, $false, "78-47=31", 2) | Out-Null
$d.Content.Find.Execute("86+4=90", $true, $false, $false, $false, $false, $true, 1, $false, "38+19=57", 2) | Out-Null
$d.Content.Find.Execute("61-32=29", $true, $false, $false, $false, $false, $true, 1, $false, "17+63=80", 2) | Out-Null
$d.Content.Find.Execute("21+77=98", $true, $false, $false, $false, $false, $true, 1, $false, "75+23=98", 2) | Out-Null
$d.Content.Find.Execute("84-82=2", $true, $false, $false, $false, $false, $true, 1, $false, "43-9=34", 2) | Out-Null
$d.Content.Find.Execute("27+38=65", $true, $false, $false, $false, $false, $true, 1, $false, "53-7=46", 2) | Out-Null
$d.Content.Find.Execute("56-18=38", $true, $false, $false, $false, $false, $true, 1, $false, "50+28=78", 2) | Out-Null
$d.Content.Find.Execute("55+14=69", $true, $false, $false, $false, $false, $true, 1, $false, "67-62=5", 2) | Out-Null
$d.Content.Find.Execute("59-35=24", $true, $false, $false, $false, $false, $true, 1, $false, "82-76=6", 2) | Out-Null
$d.Content.Find.Execute("97-5=92", $true, $false, $false, $false, $false, $true, 1, $false, "6+69=75", 2) | Out-Null
$d.Content.Find.Execute("9+51=60", $true, $false, $false, $false, $false, $true, 1, $false, "45-17=28", 2) | Out-Null
$d.Content.Find.Execute("52-11=41", $true, $false, $false, $false, $false, $true, 1, $false, "85-52=33", 2) | Out-Null
$d.Content.Find.Execute("58-58=0", $true, $false, $false, $false, $false, $true, 1, $false, "74+4=78", 2) | Out-Null
$d.Content.Find.Execute("2+67=69", $true, $false, $false, $false, $false, $true, 1, $false, "36-27=9", 2) | Out-Null
$d.Content.Find.Execute("38-6=32", $true, $false, $false, $false, $false, $true, 1, $false, "35-20=15", 2) | Out-Null
$d.Content.Find.Execute("11+46=57", $true, $false, $false, $false, $false, $true, 1, $false, "71+4=75", 2) | Out-Null
$d.Content.Find.Execute("92-43=49", $true, $false, $false, $false, $false, $true, 1, $false, "7+75=82", 2) | Out-Null
$d.Content.Find.Execute("9+68=77", $true, $false, $false, $false, $false, $true, 1, $false, "38+49=87", 2) | Out-Null
$d.Content.Find.Execute("18+76=94", $true, $false, $false, $false, $false, $true, 1, $false, "89-63=26", 2) | Out-Null
$d.Content.Find.Execute("28+27=55", $true, $false, $false, $false, $false, $true, 1, $false, "79-51=28", 2) | Out-Null
$d.Content.Find.Execute("18+56=74", $true, $false, $false, $false, $false, $true, 1, $false, "14+10=24", 2) | Out-Null
$d.Content.Find.Execute("84+6=90", $true, $false, $false, $false, $false, $true, 1, $false, "80-73=7", 2) | Out-Null
$d.Content.Find.Execute("6+43=49", $true, $false, $false, $false, $false, $true, 1, $false, "83-49=34", 2) | Out-Null
$d.Content.Find.Execute("36+50=86", $true, $false, $false, $false, $false, $true, 1, $false, "54+26=80", 2) | Out-Null
$d.Content.Find.Execute("96-29=67", $true, $false, $false, $false, $false, $true, 1, $false, "46-11=35", 2) | Out-Null
